$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (column C) was refreshed for every existing announcement
#    row (2-133): 45192 (2023-09-23) -> 45202 (2023-10-03).
$ws.Range("C2:C133").Value = 45202

# 2) Row 133 is no longer the last row, so it regains the standard 15pt
#    row height that every other data row carries.
$ws.Rows.Item(133).RowHeight = 15

# 3) Two brand new announcements were appended at the bottom of the sheet.
function Add-Announcement($Row, $Beteckning, $Datum, $Forandrad, $Lan, $Kommun, $Area, $SetHeight) {
    $ws.Cells.Item($Row, 1).Value = $Beteckning
    $ws.Cells.Item($Row, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($Row, 2).Value = $Datum
    $ws.Cells.Item($Row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($Row, 3).Value = $Forandrad
    $ws.Cells.Item($Row, 4).Value = $Lan
    $ws.Cells.Item($Row, 5).Value = $Kommun
    $ws.Cells.Item($Row, 7).Value = $Area
    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($Row, $col).Value = 0
    }
    # Column R ("Artnamn") always carries the wrap-text style, even blank.
    $ws.Cells.Item($Row, 18).WrapText = $true

    if ($SetHeight) {
        $ws.Rows.Item($Row).RowHeight = 15
    }
}

Add-Announcement 134 "A 45499-2023" 45194 45202 "STOCKHOLMS LÄN" "VALLENTUNA" 2.5 $true
Add-Announcement 135 "A 45496-2023" 45194 45202 "STOCKHOLMS LÄN" "VALLENTUNA" 1.2 $false
